$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers (row 1) ---
$ws.Range("N1").Value = "Hybrid CNOTs"
$ws.Range("O1").Value = "Hybrid Depth"

# --- Row 2: enter the formulas individually (not shared) ---
$ws.Range("N2").Formula = "=MIN(B2,E2,H2)"
$ws.Range("O2").FormulaArray = "=INDEX(CHOOSE({1,2,3},C2,F2,I2), MATCH(N2, CHOOSE({1,2,3},B2,E2,H2), 0))"

# --- Rows 3:55 - fill N as one shared-formula block (mirrors a fill-down from N3) ---
$ws.Range("N3:N55").Formula = "=MIN(B3,E3,H3)"

# --- Rows 3:55 - O column stays an individual array formula per row (legacy CSE pattern, like O2) ---
for ($r = 3; $r -le 55; $r++) {
    $formula = "=INDEX(CHOOSE({1,2,3},C$r,F$r,I$r), MATCH(N$r, CHOOSE({1,2,3},B$r,E$r,H$r), 0))"
    $ws.Range("O$r").FormulaArray = $formula
}

# --- Row 56 (mean row): extend the shared AVERAGE formula to the two new columns ---
$ws.Range("N56:O56").Formula = "=AVERAGE(N2:N55)"
$ws.Range("O56").Formula = "=AVERAGE(O2:O55)"

# --- Row 57 (improvement row): new comparison formulas ---
$ws.Range("N57").Formula = "=(`$B`$56 - N56)/`$B`$56"
$ws.Range("O57").Formula = "=(`$C`$56 - O56)/`$C`$56"

# --- View state: scroll position + selection (best effort) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 35
$ws.Range("P59").Select() | Out-Null
